# Updated via Streamlit Approval System
# Appends two new pending-approval rows (18 and 19) to the sheet, mirroring
# the existing "WGG 02" / Western Interior Designers & Marine Contractors
# rows already present in the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{
        Row = 18
        L   = "d17d0f0d-8c3f-4100-b82e-3e78cfc58db5"
        V   = 500
        X   = "PAYMENT TESTING RPA_UNIQUE_ID : 3c279024-8221-4c7d-a240-b244d4941a47"
    },
    @{
        Row = 19
        L   = "e9a48645-cfe1-4cf1-bbe6-4e38704b86a2"
        V   = 1500
        X   = "PAYMENT TESTING RPA_UNIQUE_ID : fb4f52ca-dcc2-473e-b074-af7c689919e5"
    }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value  = "WGG 02"
    $ws.Cells.Item($row, 2).Value  = "Western Interior Designers & Marine Contractors"
    $ws.Cells.Item($row, 3).Value  = "20-01-2026"
    $ws.Cells.Item($row, 4).Value  = 286962
    $ws.Cells.Item($row, 5).Value  = "Western Interior Designers & Marine Contractors"
    $ws.Cells.Item($row, 6).Value  = 34400000000
    $ws.Cells.Item($row, 7).Value  = "NEFT"
    $ws.Cells.Item($row, 8).Value  = "SBIN0003229"
    $ws.Cells.Item($row, 9).Value  = "AAAFW8862C"
    $ws.Cells.Item($row, 10).Value = "32AAAFW8862C1Z9"
    $ws.Cells.Item($row, 11).Value = ""
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = ""
    $ws.Cells.Item($row, 14).Value = ""
    $ws.Cells.Item($row, 15).Value = ""
    $ws.Cells.Item($row, 16).Value = ""
    $ws.Cells.Item($row, 17).Value = ""
    $ws.Cells.Item($row, 18).Value = ""
    $ws.Cells.Item($row, 19).Value = ""
    $ws.Cells.Item($row, 20).Value = ""
    $ws.Cells.Item($row, 21).Value = "pending"
    $ws.Cells.Item($row, 22).Value = $r.V
    $ws.Cells.Item($row, 23).Value = ""
    $ws.Cells.Item($row, 24).Value = $r.X
    $ws.Cells.Item($row, 25).Value = "HO"
    $ws.Cells.Item($row, 26).Value = 0
    $ws.Cells.Item($row, 27).Value = "midhuncraju12@gmail.com"
    $ws.Cells.Item($row, 28).Value = "ESTIMATION NOT MATCHED"
    $ws.Cells.Item($row, 29).Value = 0
    $ws.Cells.Item($row, 30).Value = 0
    $ws.Cells.Item($row, 31).Value = 0
    $ws.Cells.Item($row, 32).Value = ""
    $ws.Cells.Item($row, 33).Value = ""
    $ws.Cells.Item($row, 34).Value = ""
    $ws.Cells.Item($row, 35).Value = ""
    $ws.Cells.Item($row, 36).Value = ""
    $ws.Cells.Item($row, 37).Value = ""
    $ws.Cells.Item($row, 38).Value = ""
    $ws.Cells.Item($row, 39).Value = ""
    $ws.Cells.Item($row, 40).Value = ""
    $ws.Cells.Item($row, 41).Value = ""
}
